$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Semestre ideal" value from EF-6 to EF-4 (row 9, columns B and C)
$ws.Range("B9").Value = "EF-4"
$ws.Range("C9").Value = "EF-4"

# Delete the last row (row 24), which held the "LOM3257 - Mecânica Clássica (Requisito)" entry
$ws.Rows.Item(24).Delete()
